$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two fresh rows above the old row 230 so the existing
#    "[\InputType" thumbnail rows (old 230/231) shift down to 232/233,
#    keeping their original formatting intact.
# ---------------------------------------------------------------------------
$ws.Rows.Item(230).Insert()
$ws.Rows.Item(230).Insert()

# ---------------------------------------------------------------------------
# 2. Row 227 - "simple loop for assets" (new simplified ADAM example)
# ---------------------------------------------------------------------------
$ws.Rows.Item(227).RowHeight = 45
$ws.Range("B227").Value = "hyperlink-library"
$ws.Range("C227").Value = "simple loop for assets"

$e227 = @'
@foreach(var ${3:pic} in Adam(${1:var}, "${2:prop}").Files){
 <span>@${3:pic}.Url, @${3:pic}.FileName </span>
}
'@
$e227 = "'" + $e227
$ws.Range("E227").Value = $e227
$ws.Range("E227").WrapText = $true

$ws.Range("F227").Value = "Adam: simple example with looping ADAM assets"

# ---------------------------------------------------------------------------
# 3. Row 228 - "loop with metadata assets" (former large ADAM example, with
#    placeholders renamed from Content/Screenshots to var/prop)
# ---------------------------------------------------------------------------
$ws.Rows.Item(228).RowHeight = 150
$ws.Range("B228").Value = "hyperlink-library"
$ws.Range("C228").Value = "loop with metadata assets"

$e228 = @'
@foreach(var ${3:pic} in Adam(${1:var}, "${2:prop}").Files){
 <div style="clear: both">
  <img src="@${3:pic}.Url?w=200&h=200&mode=crop" title="@${3:pic}.FileName" style="float: right">
  <h3>@${3:pic}.Metadata.${10:Title}</h3>
  Has Meta: @${3:pic}.HasMetadata 
  <div>Description: @Html.Raw(${3:pic}.Metadata.${11:Description})</div>
 </div>
}
'@
$e228 = "'" + $e228
$ws.Range("E228").Value = $e228
$ws.Range("E228").WrapText = $true

$ws.Range("F228").Value = "Adam: Large example with looping ADAM assets"

# ---------------------------------------------------------------------------
# 4. Row 229 - "loop with type filter" (brand-new snippet, no F column)
# ---------------------------------------------------------------------------
$ws.Rows.Item(229).RowHeight = 75
$ws.Range("B229").Value = "hyperlink-library"
$ws.Range("C229").Value = "loop with type filter"

$e229 = @'
@foreach(var ${3:pic} in (AsAdam(${1:var}, "${2:prop}").Files as IEnumerable<ToSic.SexyContent.Adam.AdamFile>).Where(f => f.Type == "${4:image}")){
 <span>@${3:pic}.Url, @${3:pic}.FileName </span>
}
'@
$e229 = "'" + $e229
$ws.Range("E229").Value = $e229
$ws.Range("E229").WrapText = $true

$ws.Range("F229").ClearContents()

# ---------------------------------------------------------------------------
# 5. New rows 230/231 - "@\InputType" clones of the thumbnail-url /
#    thumbnail-IMG-tag snippets (mirrors the pre-existing "[\InputType"
#    rows that got pushed down to 232/233).
# ---------------------------------------------------------------------------
$ws.Range("A230").Value = "'@\InputType"
$ws.Range("B230").Value = "hyperlink-default"
$ws.Range("C230").Value = "thumbnail url"
$ws.Range("E230").Value = "'@`$101{var}.`$102{prop}?w=`${1:200}&h=`${2:200}&mode=`${3:crop}"
$ws.Range("F230").Value = "Thumbnail URL with crop-mode"

$ws.Range("A231").Value = "'@\InputType"
$ws.Range("B231").Value = "hyperlink-default"
$ws.Range("C231").Value = "thumbnail IMG tag"
$ws.Range("E231").Value = '<img src="@$101{var}.$102{prop}?w=${1:200}&h=${2:200}&mode=${3:crop}">'
$ws.Range("F231").Value = "Thumbnail IMG tag with crop-mode"

# ---------------------------------------------------------------------------
# 6. Grow the Table1 list object so its range/autoFilter cover the new rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F233"))

# ---------------------------------------------------------------------------
# 7. Restore the view/selection the author ended up with.
# ---------------------------------------------------------------------------
try { $excel.Goto($ws.Range("A218"), $true) } catch {}
$ws.Range("C230").Select()
